$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header, text cells) ---
# D1 is a brand-new cell: first clone C1's formatting (bold/centered/bordered
# style) onto it so it matches the rest of the header row.
$ws.Range("C1").Copy($ws.Range("D1"))

# The destination cells already carry the header style (s="1"); Excel's
# normal Value assignment would auto-convert numeric-looking strings like
# "2" into actual numbers, so stage each literal as a text formula in a
# scratch cell and paste just the resulting value back in - that keeps the
# cell's own formatting untouched while forcing a literal text value (not a
# formula) into the cell.
$ws.Range("ZZ1").Formula = '="2"'
$ws.Range("ZZ1").Copy()
$ws.Range("A1").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("ZZ1").Formula = '="1"'
$ws.Range("ZZ1").Copy()
$ws.Range("B1").PasteSpecial(-4163)

$ws.Range("ZZ1").Formula = '="0.1"'
$ws.Range("ZZ1").Copy()
$ws.Range("D1").PasteSpecial(-4163)

$ws.Range("ZZ1").ClearContents()

# --- Row 2 (numeric) ---
$ws.Range("A2").Value = 4
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0

# --- Row 3 (numeric, new row) ---
$ws.Range("A3").Value = 6
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0
